$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...NOTIFICACIONES {{NUMERO_VENDEDOR}} Y CORREO..."
#       -> "...NOTIFICACIONES AL {{NUMERO_VENDEDOR}} Y CORREO..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$found = $rng.Find.Execute(
    "NOTIFICACIONES {{NUMERO_VENDEDOR}} Y CORREO",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NOTIFICACIONES AL {{NUMERO_VENDEDOR}} Y CORREO", 2)
Write-Host "Step 1 (AL NUMERO_VENDEDOR): $found"

# ---------------------------------------------------------------------------
# 2) Drop the curly quotes around RELEVANTES:
#    QUEDANDO ESTABLECIDO QUE EN ASUNTOS SEÑALADOS "RELEVANTES", DONDE...
#       -> QUEDANDO ESTABLECIDO QUE EN ASUNTOS SEÑALADOS RELEVANTES, DONDE...
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$found = $rng.Find.Execute(
    "SEÑALADOS “RELEVANTES”, DONDE",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SEÑALADOS RELEVANTES, DONDE", 2)
Write-Host "Step 2 (drop curly quotes): $found"

# ---------------------------------------------------------------------------
# 3) Signature block: "PROMITENTE COMPRADOR" -> "PROMITENTE {{SEXO_4}}"
#                      "PROMITENTE VENDEDOR"  -> "PROMITENTE {{SEXO_2}}"
#    (only the signature-block pair, not the earlier body mention of
#    "PROMITENTE VENDEDOR" in the custody clause). Replace just the
#    trailing word so the leading "PROMITENTE" run/its formatting carries
#    over untouched, same as the source edit did.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.Text = " COMPRADOR"
$found = $rng.Find.Execute()
Write-Host "Step 3a (locate PROMITENTE COMPRADOR): $found"
if ($found) {
    $rng.Text = " {{SEXO_4}}"

    $rng.Collapse(0)
    $rng.Find.ClearFormatting()
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 0
    $rng.Find.Text = " VENDEDOR"
    $found2 = $rng.Find.Execute()
    Write-Host "Step 3b (locate PROMITENTE VENDEDOR after it): $found2"
    if ($found2) {
        $rng.Text = " {{SEXO_2}}"
    }
}
